$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("N52").ClearContents()
$ws.Range("H53").Value = 2326.3333
$ws.Range("I53").Value = 1881.0834
$ws.Range("J53").Value = 3216.8333
$ws.Range("K53").Value = 1881.0834
$ws.Range("L53").Value = 3216.8333
$ws.Range("M53").Value = -1244.0834
$ws.Range("N53").Value = -4490.8333
$ws.Range("H55").Value = 906.1875
$ws.Range("I55").Value = 326.875
$ws.Range("J55").Value = 1485.5
$ws.Range("K55").Value = 326.875
$ws.Range("L55").Value = 1485.5
$ws.Range("M55").Value = -112.875
$ws.Range("N55").Value = -1913.5
$ws.Range("H112").Value = 2063.56
$ws.Range("I112").Value = 2137.6
$ws.Range("K112").Value = 6412.799999999999
$ws.Range("M112").Value = -5304.799999999999
$ws.Range("H132").Value = 41799.24
$ws.Range("I132").Value = 51092.3
$ws.Range("J132").Value = 4627
$ws.Range("K132").Value = 153276.9
$ws.Range("L132").Value = 13881
$ws.Range("M132").Value = -150746.9
$ws.Range("N132").Value = -18941
$ws.Range("H135").Value = 434.53333
$ws.Range("I135").Value = 434.53333
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 3910.79997
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -1375.79997
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 2216.5
$ws.Range("I137").Value = 2216.5
$ws.Range("K137").Value = 6649.5
$ws.Range("M137").Value = -4099.5
$ws.Range("H138").Value = 2375.3381
$ws.Range("I138").Value = 909.75
$ws.Range("J138").Value = 2826.2886
$ws.Range("K138").Value = 2729.25
$ws.Range("L138").Value = 8478.8658
$ws.Range("M138").Value = 2410.75
$ws.Range("N138").Value = -18758.8658

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1632.0222
$ws.Range("I32").Value = 1632.0222
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1632.0222
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1345.0222
$ws.Range("N32").ClearContents()
$ws.Range("H61").Value = 3264.647
$ws.Range("I61").Value = 3149.3
$ws.Range("J61").Value = 4129.75
$ws.Range("K61").Value = 3149.3
$ws.Range("L61").Value = 4129.75
$ws.Range("M61").Value = -2937.3
$ws.Range("N61").Value = -4553.75
$ws.Range("H136").Value = 3264.647
$ws.Range("I136").Value = 3149.3
$ws.Range("J136").Value = 4129.75
$ws.Range("K136").Value = 9447.900000000001
$ws.Range("L136").Value = 12389.25
$ws.Range("M136").Value = -6897.900000000001
$ws.Range("N136").Value = -17489.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1939.8572
$ws.Range("I107").Value = 642.93335
$ws.Range("J107").Value = 5182.1665
$ws.Range("K107").Value = 642.93335
$ws.Range("L107").Value = 5182.1665
$ws.Range("M107").Value = 1277.06665
$ws.Range("N107").Value = -9022.166499999999
$ws.Range("H134").Value = 58335736
$ws.Range("I134").Value = 35717144
$ws.Range("J134").Value = 111112450
$ws.Range("K134").Value = 107151432
$ws.Range("L134").Value = 333337350
$ws.Range("M134").Value = -107148897
$ws.Range("N134").Value = -333342420

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 965.5833
$ws.Range("I16").Value = 944.36365
$ws.Range("K16").Value = 944.36365
$ws.Range("M16").Value = -657.36365
$ws.Range("H31").Value = 2916.4167
$ws.Range("I31").Value = 2559.6
$ws.Range("J31").Value = 3171.2856
$ws.Range("K31").Value = 2559.6
$ws.Range("L31").Value = 3171.2856
$ws.Range("M31").Value = -2264.6
$ws.Range("N31").Value = -3761.2856
$ws.Range("H34").Value = 2916.4167
$ws.Range("I34").Value = 2559.6
$ws.Range("J34").Value = 3171.2856
$ws.Range("K34").Value = 2559.6
$ws.Range("L34").Value = 3171.2856
$ws.Range("M34").Value = -2357.6
$ws.Range("N34").Value = -3575.2856
$ws.Range("H58").Value = 1743.7949
$ws.Range("I58").Value = 1466.9395
$ws.Range("J58").Value = 3266.5
$ws.Range("K58").Value = 1466.9395
$ws.Range("L58").Value = 3266.5
$ws.Range("M58").Value = -1263.9395
$ws.Range("N58").Value = -3672.5
$ws.Range("H113").Value = 965.5833
$ws.Range("I113").Value = 944.36365
$ws.Range("K113").Value = 944.36365
$ws.Range("M113").Value = 1225.63635
$ws.Range("H134").Value = 2085190.1
$ws.Range("I134").Value = 1715.2572
$ws.Range("J134").Value = 7694545.5
$ws.Range("K134").Value = 5145.7716
$ws.Range("L134").Value = 23083636.5
$ws.Range("M134").Value = -2610.7716
$ws.Range("N134").Value = -23088706.5
$ws.Range("H136").Value = 1743.7949
$ws.Range("I136").Value = 1466.9395
$ws.Range("J136").Value = 3266.5
$ws.Range("K136").Value = 4400.818499999999
$ws.Range("L136").Value = 9799.5
$ws.Range("M136").Value = -1850.818499999999
$ws.Range("N136").Value = -14899.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 6862.909
$ws.Range("I81").Value = 4499.2856
$ws.Range("J81").Value = 10999.25
$ws.Range("K81").Value = 13497.8568
$ws.Range("L81").Value = 32997.75
$ws.Range("M81").Value = -12374.8568
$ws.Range("N81").Value = -35243.75
$ws.Range("H84").Value = 6862.909
$ws.Range("I84").Value = 4499.2856
$ws.Range("J84").Value = 10999.25
$ws.Range("K84").Value = 40493.5704
$ws.Range("L84").Value = 98993.25
$ws.Range("M84").Value = -34877.5704
$ws.Range("N84").Value = -110225.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2300
$ws.Range("I102").Value = 2300
$ws.Range("K102").Value = 2300
$ws.Range("M102").Value = -678
$ws.Range("H132").Value = 2275.0227
$ws.Range("I132").Value = 2342.7
$ws.Range("J132").Value = 1598.25
$ws.Range("K132").Value = 7028.099999999999
$ws.Range("L132").Value = 4794.75
$ws.Range("M132").Value = -4498.099999999999
$ws.Range("N132").Value = -9854.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7614.3335
$ws.Range("I40").Value = 6338.75
$ws.Range("J40").Value = 8634.799999999999
$ws.Range("K40").Value = 6338.75
$ws.Range("L40").Value = 8634.799999999999
$ws.Range("M40").Value = -6202.75
$ws.Range("N40").Value = -8906.799999999999
$ws.Range("H82").Value = 2033.25
$ws.Range("I82").Value = 1189.6
$ws.Range("J82").Value = 6251.5
$ws.Range("K82").Value = 1189.6
$ws.Range("L82").Value = 6251.5
$ws.Range("M82").Value = -828.5999999999999
$ws.Range("N82").Value = -6973.5
$ws.Range("H85").Value = 2033.25
$ws.Range("I85").Value = 1189.6
$ws.Range("J85").Value = 6251.5
$ws.Range("K85").Value = 1189.6
$ws.Range("L85").Value = 6251.5
$ws.Range("M85").Value = 58.40000000000009
$ws.Range("N85").Value = -8747.5
$ws.Range("H136").Value = 26319020
$ws.Range("I136").Value = 2960.742
$ws.Range("J136").Value = 142861570
$ws.Range("K136").Value = 8882.226000000001
$ws.Range("L136").Value = 428584710
$ws.Range("M136").Value = -6332.226000000001
$ws.Range("N136").Value = -428589810

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 13476
$ws.Range("J45").Value = 13476
$ws.Range("L45").Value = 13476
$ws.Range("N45").Value = -14458
$ws.Range("H107").Value = 1526.1364
$ws.Range("I107").Value = 1465.2222
$ws.Range("J107").Value = 1800.25
$ws.Range("K107").Value = 4395.6666
$ws.Range("L107").Value = 5400.75
$ws.Range("M107").Value = -2475.6666
$ws.Range("N107").Value = -9240.75
$ws.Range("H113").Value = 3100.4
$ws.Range("I113").Value = 2834
$ws.Range("J113").Value = 3500
$ws.Range("K113").Value = 8502
$ws.Range("L113").Value = 10500
$ws.Range("M113").Value = -6332
$ws.Range("N113").Value = -14840
$ws.Range("H135").Value = 50683.625
$ws.Range("J135").Value = 50683.625
$ws.Range("L135").Value = 50683.625
$ws.Range("N135").Value = -60823.625
